# Populate the "particle throughput" sample table on Sheet1.
# Columns: A = Particle Count, B = Peak Performance (MegaParticles/Sec)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (written as shared strings by the engine)
$ws.Range("A1").Value = "Particle Count"
$ws.Range("B1").Value = "Peak Performance (MegaParticles/Sec)"

# Data rows
$ws.Range("A2").Value = 1048576
$ws.Range("B2").Value = 1.3

$ws.Range("A3").Value = 524288
$ws.Range("B3").Value = 1.2

$ws.Range("A4").Value = 262144
$ws.Range("B4").Value = 0.9

$ws.Range("A5").Value = 131072
$ws.Range("B5").Value = 0.7

# Particle counts get thousands-separator formatting (builtin numFmtId 3)
$ws.Range("A2:A5").NumberFormat = "#,##0"

# Widen the two columns to fit their (now longer) header text
$ws.Columns("A:B").AutoFit()

# Leftover cursor position from the editing session
[void]$ws.Range("G6").Select()
